$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Footers(1) -> footer2.xml (default footer), Footers(2) -> footer1.xml (first-page footer).
# Both hold the Pearson Edexcel logo; rename its picture from image2.png to image1.png.
$f1 = $sec.Footers(1)
if ($f1.Exists -and $f1.Range.InlineShapes.Count -ge 1) {
  $f1.Range.InlineShapes(1).Name = "image1.png"
}

$f2 = $sec.Footers(2)
if ($f2.Exists -and $f2.Range.InlineShapes.Count -ge 1) {
  $f2.Range.InlineShapes(1).Name = "image1.png"
}

# Headers(2) -> header1.xml (first-page header) holds the BTec logo; rename its
# picture from image1.jpg to image2.jpg.
$h2 = $sec.Headers(2)
if ($h2.Exists -and $h2.Range.InlineShapes.Count -ge 1) {
  $h2.Range.InlineShapes(1).Name = "image2.jpg"
}
